$d = $word.ActiveDocument

$replacements = @(
    @{old="2025-03-24 Monday"; new="2025-03-25 Tuesday"},
    @{old="421÷7="; new="747÷7="},
    @{old="536÷7="; new="224÷8="},
    @{old="699÷5="; new="797÷9="},
    @{old="412÷4="; new="390÷4="},
    @{old="115÷6="; new="853÷4="},
    @{old="465÷9="; new="144÷3="},
    @{old="999÷8="; new="291÷9="},
    @{old="648÷7="; new="218÷9="},
    @{old="367÷7="; new="345÷6="},
    @{old="669÷3="; new="380÷3="},
    @{old="128÷4="; new="635÷6="},
    @{old="714÷9="; new="147÷6="},
    @{old="977÷2="; new="303÷5="},
    @{old="919÷6="; new="756÷9="},
    @{old="980÷4="; new="285÷8="},
    @{old="644÷8="; new="691÷2="},
    @{old="319÷4="; new="214÷7="},
    @{old="367÷5="; new="780÷8="},
    @{old="729÷7="; new="536÷2="},
    @{old="909÷6="; new="542÷4="},
    @{old="894÷7="; new="144÷5="},
    @{old="182÷6="; new="778÷7="},
    @{old="392÷4="; new="894÷2="},
    @{old="827÷8="; new="506÷2="},
    @{old="655÷3="; new="933÷2="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
